# Auto-generated Excel COM-interop script applying the Siren_Profits.xlsx diff
# to the corresponding sheets in this workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 51
$ws.Range("H51").Value = 57945
$ws.Range("I51").Value = 2535.25
$ws.Range("J51").Value = 102272.8
$ws.Range("K51").Value = 2535.25
$ws.Range("L51").Value = 102272.8
$ws.Range("M51").Value = -2051.25
$ws.Range("N51").Value = -103240.8
# row 69
$ws.Range("H69").Value = 18999.2
$ws.Range("J69").Value = 39999.5
$ws.Range("L69").Value = 119998.5
$ws.Range("N69").Value = -121746.5
# row 72
$ws.Range("H72").Value = 18999.2
$ws.Range("J72").Value = 39999.5
$ws.Range("L72").Value = 359995.5
$ws.Range("N72").Value = -368731.5
# row 99
$ws.Range("H99").Value = 1823605.5
$ws.Range("J99").Value = 1545
$ws.Range("L99").Value = 4635
$ws.Range("N99").Value = -7631
# row 112
$ws.Range("H112").Value = 3800
$ws.Range("J112").Value = 2500
$ws.Range("L112").Value = 7500
$ws.Range("N112").Value = -9716
# row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
# row 134
$ws.Range("H134").Value = 92985.19
$ws.Range("J134").Value = 92985.19
$ws.Range("L134").Value = 92985.19
$ws.Range("N134").Value = -103125.19
# row 137
$ws.Range("H137").Value = 752569.9399999999
$ws.Range("I137").Value = 1022107
$ws.Range("K137").Value = 3066321
$ws.Range("M137").Value = -3063771

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3544.08
$ws.Range("I32").Value = 3544.08
$ws.Range("K32").Value = 3544.08
$ws.Range("M32").Value = -3257.08
# row 61
$ws.Range("H61").Value = 5381.7334
$ws.Range("I61").Value = 5831.2856
$ws.Range("K61").Value = 5831.2856
$ws.Range("M61").Value = -5619.2856
# row 64
$ws.Range("H64").Value = 45000
$ws.Range("I64").Value = 45000
$ws.Range("K64").Value = 45000
$ws.Range("M64").Value = -44752
# row 67
$ws.Range("H67").Value = 45000
$ws.Range("I67").Value = 45000
$ws.Range("K67").Value = 45000
$ws.Range("M67").Value = -44142
# row 74
$ws.Range("H74").Value = 3747.8306
$ws.Range("I74").Value = 15915.125
$ws.Range("K74").Value = 15915.125
$ws.Range("M74").Value = -15041.125
# row 77
$ws.Range("H77").Value = 3747.8306
$ws.Range("I77").Value = 15915.125
$ws.Range("K77").Value = 79575.625
$ws.Range("M77").Value = -75207.625
# row 97
$ws.Range("H97").Value = 9095711
$ws.Range("I97").Value = 5030.6665
$ws.Range("K97").Value = 5030.6665
$ws.Range("M97").Value = -4534.6665
# row 135
$ws.Range("H135").Value = 197499
$ws.Range("J135").Value = 197499
$ws.Range("L135").Value = 197499
$ws.Range("N135").Value = -207639
# row 136
$ws.Range("H136").Value = 5381.7334
$ws.Range("I136").Value = 5831.2856
$ws.Range("K136").Value = 17493.8568
$ws.Range("M136").Value = -14943.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 106
$ws.Range("H106").Value = 65916.5
$ws.Range("J106").Value = 65916.5
$ws.Range("L106").Value = 65916.5
$ws.Range("N106").Value = -68440.5
# row 107
$ws.Range("H107").Value = 1389.7778
$ws.Range("I107").Value = 1438.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1438.5
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 481.5
$ws.Range("N107").Value = -4840
# row 140
$ws.Range("H140").Value = 122945
$ws.Range("I140").Value = 100000
$ws.Range("K140").Value = 100000
$ws.Range("M140").Value = -94820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 4879.3
$ws.Range("I31").Value = 1650
$ws.Range("J31").Value = 5686.625
$ws.Range("K31").Value = 1650
$ws.Range("L31").Value = 5686.625
$ws.Range("M31").Value = -1355
$ws.Range("N31").Value = -6276.625
# row 34
$ws.Range("H34").Value = 4879.3
$ws.Range("I34").Value = 1650
$ws.Range("J34").Value = 5686.625
$ws.Range("K34").Value = 1650
$ws.Range("L34").Value = 5686.625
$ws.Range("M34").Value = -1448
$ws.Range("N34").Value = -6090.625
# row 95
$ws.Range("H95").Value = 312015000
$ws.Range("J95").Value = 312015000
$ws.Range("L95").Value = 312015000
$ws.Range("N95").Value = -312020492
# row 119
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# row 135
$ws.Range("H135").Value = 57145.5
$ws.Range("J135").Value = 57145.5
$ws.Range("L135").Value = 57145.5
$ws.Range("N135").Value = -67285.5
# row 141
$ws.Range("H141").Value = 411739
$ws.Range("J141").Value = 517931.1
$ws.Range("L141").Value = 517931.1
$ws.Range("N141").Value = -528291.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 131
$ws.Range("H131").Value = 5012.4585
$ws.Range("I131").Value = 8269.454
$ws.Range("K131").Value = 24808.362
$ws.Range("M131").Value = -19768.362

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 13489.333
$ws.Range("I70").Value = 8399.4
$ws.Range("K70").Value = 8399.4
$ws.Range("M70").Value = -8129.4
# row 73
$ws.Range("H73").Value = 13489.333
$ws.Range("I73").Value = 8399.4
$ws.Range("K73").Value = 8399.4
$ws.Range("M73").Value = -7463.4
# row 123
$ws.Range("H123").Value = 46000
$ws.Range("J123").Value = 46000
$ws.Range("L123").Value = 46000
$ws.Range("N123").Value = -50900
# row 126
$ws.Range("H126").Value = 28915.666
$ws.Range("I126").Value = 61499.5
$ws.Range("J126").Value = 12623.75
$ws.Range("K126").Value = 184498.5
$ws.Range("L126").Value = 37871.25
$ws.Range("M126").Value = -182028.5
$ws.Range("N126").Value = -42811.25
# row 132
$ws.Range("H132").Value = 3492.5217
$ws.Range("I132").Value = 3544.25
$ws.Range("K132").Value = 10632.75
$ws.Range("M132").Value = -8102.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 93
$ws.Range("H93").Value = 3775.6216
$ws.Range("I93").Value = 4611.68
$ws.Range("J93").Value = 2033.8334
$ws.Range("K93").Value = 4611.68
$ws.Range("L93").Value = 2033.8334
$ws.Range("M93").Value = -3363.68
$ws.Range("N93").Value = -4529.8334
# row 132
$ws.Range("H132").Value = 407375.28
$ws.Range("I132").Value = 936663.5600000001
$ws.Range("J132").Value = 4108
$ws.Range("K132").Value = 2809990.68
$ws.Range("L132").Value = 12324
$ws.Range("M132").Value = -2807460.68
$ws.Range("N132").Value = -17384
# row 136
$ws.Range("H136").Value = 8423.24
$ws.Range("I136").Value = 5871.933
$ws.Range("K136").Value = 17615.799
$ws.Range("M136").Value = -15065.799

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 46
$ws.Range("H46").Value = 43124.5
$ws.Range("J46").Value = 26249
$ws.Range("L46").Value = 26249
$ws.Range("N46").Value = -26711
# row 132
$ws.Range("H132").Value = 7891.763
$ws.Range("I132").Value = 9545.386
$ws.Range("K132").Value = 28636.158
$ws.Range("M132").Value = -26106.158
# row 134
$ws.Range("H134").Value = 43124.5
$ws.Range("J134").Value = 26249
$ws.Range("L134").Value = 78747
$ws.Range("N134").Value = -83817
# row 136
$ws.Range("H136").Value = 1102139
$ws.Range("I136").Value = 1284912.1
$ws.Range("K136").Value = 3854736.3
$ws.Range("M136").Value = -3852186.3
